$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# Insert a new row before row 18 (shifts rows 18-69 down to 19-70)
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with "Router ID", matching the style/format of the
# row that used to be there (which is now row 19, "VM name")
$ws.Cells.Item(18, 1).Value = "Router ID"

# Apply the same style as the row below it (row 19, which had the style
# previously used by row 18 before the insert) so formatting matches.
$ws.Cells.Item(18, 1).Style = $ws.Cells.Item(19, 1).Style
$ws.Cells.Item(18, 2).Style = $ws.Cells.Item(19, 2).Style
$ws.Cells.Item(18, 3).Style = $ws.Cells.Item(19, 3).Style

# Add the comment for the new Router ID cell
$comment = $ws.Cells.Item(18, 1).AddComment("Required IPv4 address when using an IPv6 system IP address [default: (system_ip)]")
